# Fix for Excel files with empty rows ("fix para exceles con rows vacias")
#
# Adds a few whitespace-only placeholder values below the existing data
# (so rows that were previously completely empty now carry a value) and
# marks a further cell (C7) with an underlined style, ready for manual
# entry/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: it determines the order new entries are appended to the
# shared-strings table (D3's 4-space value is interned before C3's
# 2-space value).
$ws.Range("D3").Value = "    "
$ws.Range("C3").Value = "  "
$ws.Range("C4").Value = "   "
$ws.Range("C5").Value = "   "

# Leave C7 blank but give it an underline font so it picks up a dedicated
# (new) cell style, then leave the selection parked there.
$ws.Range("C7").Font.Underline = $true
$ws.Range("C7").Select() | Out-Null
